$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.542.97'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '1.825.01'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '315.51'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '0.5114'
$ws.Range('E7').Value = '  -5.41%  '
$ws.Range('D8').Value = '0.3955'
$ws.Range('E8').Value = '  -1.55%  '
$ws.Range('D9').Value = '0.08222'
$ws.Range('E9').Value = '  +6.26%  '
$ws.Range('D10').Value = '1.116'
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').Value = '41.72'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('D12').Value = '21.19'
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').Value = '6.344'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').Value = '7.555'
$ws.Range('E15').Value = '  -1.33%  '
$ws.Range('D16').Value = '1.823.14'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '0.00001126'
$ws.Range('E17').Value = '  +3.28%  '
$ws.Range('D18').Value = '92.95'
$ws.Range('E18').Value = '  +3.13%  '
$ws.Range('D19').Value = '0.06648'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').Value = '6.097'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').Value = '28.581.04'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = '11.44'
$ws.Range('E24').Value = '  +2.03%  '
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('D26').Value = '21.38'
$ws.Range('E26').Value = '  +2.65%  '
$ws.Range('D27').Value = '156.50'
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('D28').Value = '2.034.86'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').Value = '2.419'
$ws.Range('E29').Value = '  -2.11%  '
$ws.Range('D30').Value = '126.95'
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('D31').Value = '1.114'
$ws.Range('E31').Value = '  -1.49%  '
$ws.Range('D32').Value = '0.1088'
$ws.Range('E32').Value = '  -2.53%  '
$ws.Range('D33').Value = '5.768'
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('D35').Value = '0.07050'
$ws.Range('E35').Value = '  -6.03%  '
$ws.Range('D36').Value = '0.2231'
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('D37').Value = '5.285'
$ws.Range('E37').Value = '  +1.13%  '
$ws.Range('D38').Value = '0.02356'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D39').Value = '8.810'
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('D40').Value = '0.6336'
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').Value = '11.29'
$ws.Range('E41').Value = '  -1.03%  '
$ws.Range('D42').Value = '1.183'
$ws.Range('D43').Value = '1.400'
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D44').Value = '13.55'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').Value = '0.5950'
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('E46').Value = '  +0.60%  '
$ws.Range('D47').Value = '125.14'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('D48').Value = '1.997'
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('D49').Value = '1.192'
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('D50').Value = '0.06948'
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('D51').Value = '1.084'
$ws.Range('E51').Value = '  +4.23%  '
